$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (currentMarketValue, profitLossAmount, profitLossPercentage)
$ws.Range("D2").Value = 14890378.1
$ws.Range("E2").Value = -514422.25
$ws.Range("F2").Value = -32.03

# Update row 3 values (currentMarketValue, profitLossAmount, profitLossPercentage)
$ws.Range("D3").Value = 17480718.7
$ws.Range("E3").Value = -1583190.08
$ws.Range("F3").Value = 2758.46
